# The blog template originally contained a pile of experimental/duplicate
# paragraphs exercising different placeholder syntaxes ($... / {{...}} /
# mixed with field codes). The final template should just have two simple
# paragraphs: {{TITLE}} and {{CONTENT}}, each a single plain-text run.

$d = $word.ActiveDocument

# 1. Drop every paragraph after the second one (paragraphs 3-6 in the
#    original file hold the duplicate/experimental placeholder variants).
if ($d.Paragraphs.Count -gt 2) {
    $firstExtra = $d.Paragraphs.Item(3)
    $tail = $d.Range($firstExtra.Range.Start, $d.Content.End)
    $tail.Delete()
}

# 2. Paragraph 1 currently reads "$" + "{TITLE}" across two runs -> collapse
#    it to a single run reading "{{TITLE}}".
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.MoveEnd(1, -1) | Out-Null
$r1.Text = "{{TITLE}}"

# 3. Paragraph 2 currently reads "$" + "{CONTENT}" followed by a
#    { {{TITLE}} } field-code run -> remove the field, then collapse the
#    paragraph to a single run reading "{{CONTENT}}".
$p2 = $d.Paragraphs.Item(2)
while ($p2.Range.Fields.Count -gt 0) {
    $p2.Range.Fields.Item(1).Delete()
}
$r2 = $p2.Range
$r2.MoveEnd(1, -1) | Out-Null
$r2.Text = "{{CONTENT}}"
